$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update device IDs (column D) for the ANDROID device rows
$ws.Range("D4").Value = "234c19cb26017ece"
$ws.Range("D5").Value = "234c19cb26017ece"

# Update device ID for the IOS device row
$ws.Range("D6").Value = "fd76a9d32fb7cc6eb6284cbcab936bc97dcfce35"

# Update device names (column H) for the ANDROID device rows
$ws.Range("H4").Value = "SM-G960F"
$ws.Range("H5").Value = "SM-G960F"

# Match the updated selection left by the edit
$ws.Range("E6").Select()
